# Hard coded values have been changed:
#  - "Test Cases" Runmode for TC_001 flipped from Y to N
#  - "Test Data" stray hard-coded "PASS" status cell cleared

$wb = $excel.ActiveWorkbook

# --- "Test Steps" sheet: rest the selection on A15 ---
$ws2 = $wb.Worksheets.Item("Test Steps")
[void]$ws2.Range("A15").Select()

# --- "Test Data" sheet: clear the hard-coded "PASS" value, move to D11 ---
$ws3 = $wb.Worksheets.Item("Test Data")
$ws3.Range("E3").ClearContents()
[void]$ws3.Range("D11").Select()

# --- "Test Cases" sheet: Runmode hard-coded value Y -> N ---
$ws1 = $wb.Worksheets.Item("Test Cases")
$ws1.Range("C2").Value = "N"
[void]$ws1.Range("C2").Select()
